# Update metrics across multiple worksheets, adding new run data
# (re-ran benchmarks, column width tweak on a couple sheets to fit graph)

$wb = $excel.ActiveWorkbook

function Set-Metrics {
    param(
        [string]$SheetName,
        [string[]]$Values
    )
    $ws = $wb.Worksheets.Item($SheetName)
    for ($i = 0; $i -lt $Values.Length; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 2).Value = [double]$Values[$i]
    }
}

# "Bidirectional A" sheet
Set-Metrics "Bidirectional A" @(
    "0.0001963750000868458",
    "0.0002029580009548226",
    "5431.7",
    "5431.7",
    "0.00315093994140625",
    "0.00309600830078125",
    "0.0084991455078125",
    "0.0084991455078125",
    "4.362500112620182e-05",
    "6.170900087454356e-05",
    "2.179581254040386e-05",
    "4.326136256622703e-05",
    "80"
)

# "D Lite" sheet
Set-Metrics "D Lite" @(
    "0.000353749999703723",
    "0.001694665999821154",
    "5431.7",
    "5431.7",
    "0.003082275390625",
    "0.003082275390625",
    "0.07830810546875",
    "0.07805328369140625",
    "0.0002877090009860694",
    "0.0008845000011206139",
    "0.0001083754750425214",
    "4.761206259900064e-05",
    "80"
)

# "IDA" sheet (also narrows column B from 24 to 23)
Set-Metrics "IDA" @(
    "0.0001899170001706807",
    "0.009287000000767875",
    "5431.7",
    "5431.7",
    "0.003082275390625",
    "0.003082275390625",
    "0.00305938720703125",
    "0.00300445556640625",
    "6.512499930977356e-05",
    "0.005360707998988801",
    "0.008276711987468844",
    "4.337861246312968e-05",
    "80"
)
$wsIDA = $wb.Worksheets.Item("IDA")
$wsIDA.Range("B:B").ColumnWidth = 22.166666666666668

# "SMA" sheet
Set-Metrics "SMA" @(
    "0.0001857920015027048",
    "6.320799911918584e-05",
    "5431.7",
    "6729.199999999999",
    "0.00315093994140625",
    "0.00313262939453125",
    "0.00220489501953125",
    "0.00220489501953125",
    "4.333299875725061e-05",
    "3.950000063923653e-05",
    "0.0001048161124572289",
    "6.615983757001232e-05",
    "80"
)

# "RTAA (L=25, M=3)" sheet (also narrows column B from 24 to 23)
Set-Metrics "RTAA (L=25, M=3)" @(
    "0.0004321249998611165",
    "0.0004928750004182803",
    "5431.7",
    "5431.7",
    "0.00313568115234375",
    "0.00309295654296875",
    "0.00685882568359375",
    "0.00685882568359375",
    "4.429200089361984e-05",
    "0.0001279580010304926",
    "0.0001367577998962588",
    "4.306251253183291e-05",
    "80"
)
$wsRTAA = $wb.Worksheets.Item("RTAA (L=25, M=3)")
$wsRTAA.Range("B:B").ColumnWidth = 22.166666666666668
